$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" column (D) cells being updated hold plain numeric-looking
# text (e.g. "64.015.56", "1.00", "0.110") in the source workbook - they
# are inline strings, not real numbers (note the thousands-dot style and
# trailing zeros). Force Text format first so Excel/COM does not coerce
# them into Number cells (which would also silently drop the formatting,
# e.g. "1.00" -> 1).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '64.015.56'
$ws.Range('D3').Value = '2.753.43'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '577.14'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '157.64'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').Value = '0.110'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').Value = '0.384'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').Value = '5.59'
$ws.Range('E12').Value = '  -16.54%  '
$ws.Range('D13').Value = '3.241.26'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '26.79'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '63.692.23'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '0.0000151'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').Value = '2.754.74'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '12.13'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').Value = '4.86'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '358.32'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').Value = '0.545'
$ws.Range('E22').Value = '  +2.09%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '65.66'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '8.46'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').Value = '0.0₃0904'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').Value = '1.94'
$ws.Range('E29').Value = '  -2.76%  '
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '170.72'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.21'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').Value = '20.29'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').Value = '4.92'
$ws.Range('E34').Value = '  +3.26%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').Value = '1.45'
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('D37').Value = '1.80'
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').Value = '0.985'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').Value = '6.20'
$ws.Range('E39').Value = '  +11.75%  '
$ws.Range('D40').Value = '4.15'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').Value = '326.27'
$ws.Range('E41').Value = '  -4.99%  '
$ws.Range('D42').Value = '39.29'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').Value = '21.44'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').Value = '0.0592'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Value = '21.68'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').Value = '0.0255'
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('D47').Value = '136.33'
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('D48').Value = '0.630'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').Value = '0.101'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('E51').Value = '  +0.57%  '
